$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamps = @{
    2  = "2025-11-03T00:10:42.490533"
    3  = "2025-11-03T00:10:42.490533"
    4  = "2025-11-03T00:10:42.490533"
    5  = "2025-11-03T00:10:42.491911"
    6  = "2025-11-03T00:10:42.491911"
    7  = "2025-11-03T00:10:42.491911"
    8  = "2025-11-03T00:10:42.491911"
    9  = "2025-11-03T00:10:42.491911"
    10 = "2025-11-03T00:10:42.491911"
    11 = "2025-11-03T00:10:42.492909"
    12 = "2025-11-03T00:10:42.492909"
    13 = "2025-11-03T00:10:42.492909"
    14 = "2025-11-03T00:10:42.492909"
    15 = "2025-11-03T00:10:42.492909"
    16 = "2025-11-03T00:10:42.492909"
    17 = "2025-11-03T00:10:42.492909"
    18 = "2025-11-03T00:10:42.492909"
    19 = "2025-11-03T00:10:42.492909"
    20 = "2025-11-03T00:10:42.492909"
    21 = "2025-11-03T00:10:42.493911"
    22 = "2025-11-03T00:10:42.493911"
    23 = "2025-11-03T00:10:42.493911"
    24 = "2025-11-03T00:10:42.493911"
    25 = "2025-11-03T00:10:42.493911"
    26 = "2025-11-03T00:10:42.493911"
    27 = "2025-11-03T00:10:42.493911"
    28 = "2025-11-03T00:10:42.493911"
    29 = "2025-11-03T00:10:42.494903"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}

$wb.Save()
